$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Row 5 (CAIZA COLLAGUAZO ROCIO PILAR)
$ws1.Range("H5").Value = 811.8
$ws1.Range("P5").Value = 550.63

# Row 21 (TAMAYO VILLACIS EDWIN XAVIER)
$ws1.Range("M21").Value = 160.38
$ws1.Range("P21").Value = 550.63

# Row 26 (summary counts "X de 24")
$ws1.Range("H26").Value = "1 de 24"
$ws1.Range("M26").Value = "4 de 24"
$ws1.Range("P26").Value = "2 de 24"

# --- Sheet: VENTA MENSUAL ---
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")

$ws2.Range("F5").Value = 1362.43
$ws2.Range("F21").Value = 711.01
$ws2.Range("F26").Value = 25559.02

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# Row 6 - INODOROS
$ws3.Range("D6").Value = 642.61
$ws3.Range("E6").Value = 264.556108615601
$ws3.Range("F6").Value = 0.708370819739582

# Row 8 - NO RESURTIBLES
$ws3.Range("D8").Value = 1101.26
$ws3.Range("E8").Value = -434.673172431852
$ws3.Range("F8").Value = 1.652087851807143

# Row 12 - PORCELANATO
$ws3.Range("D12").Value = 22615.46
$ws3.Range("E12").Value = 5339.52
$ws3.Range("F12").Value = 0.8089957495945266

# Row 14 - TOTAL
$ws3.Range("D14").Value = 24366.13
$ws3.Range("E14").Value = 17837.25110009469
$ws3.Range("F14").Value = 0.5773501877067695
